$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.240041375160217
$ws.Range("B1").Value = 2.345365762710571
$ws.Range("C1").Value = 3.655198574066162
$ws.Range("D1").Value = 3.490769386291504
$ws.Range("E1").Value = 1.143975257873535
